$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 25,14
$arr[0,0] = 14.734287
$arr[0,1] = 44.202861
$arr[0,2] = 0.4000023944294819
$arr[0,3] = 0.400002394429482
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 68.637375
$arr[0,7] = 205.912125
$arr[0,8] = 0.5415701538216162
$arr[0,9] = 0.5415701538216162
$arr[0,10] = 1011.322782176625
$arr[0,11] = 9101.905039589625
$arr[0,12] = 0.2166293582801893
$arr[0,13] = 0.2166293582801894
$arr[1,0] = 14.734287
$arr[1,1] = 44.202861
$arr[1,2] = 0.4000023944294819
$arr[1,3] = 0.400002394429482
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 11.05007466666667
$arr[1,7] = 33.150224
$arr[1,8] = 0.08718851262838957
$arr[1,9] = 0.08718851262838957
$arr[1,10] = 162.814971510096
$arr[1,11] = 1465.334743590864
$arr[1,12] = 0.03487561381810095
$arr[1,13] = 0.03487561381810096
$arr[2,0] = 14.734287
$arr[2,1] = 44.202861
$arr[2,2] = 0.4000023944294819
$arr[2,3] = 0.400002394429482
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 16.21089566666667
$arr[2,7] = 48.632687
$arr[2,8] = 0.1279089892319285
$arr[2,9] = 0.1279089892319285
$arr[2,10] = 238.855989279723
$arr[2,11] = 2149.703903517507
$arr[2,12] = 0.05116390196182622
$arr[2,13] = 0.05116390196182623
$arr[3,0] = 14.734287
$arr[3,1] = 44.202861
$arr[3,2] = 0.4000023944294819
$arr[3,3] = 0.400002394429482
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 20.32546233333333
$arr[3,7] = 60.976387
$arr[3,8] = 0.1603741949973873
$arr[3,9] = 0.1603741949973873
$arr[3,10] = 299.481195427023
$arr[3,11] = 2695.330758843207
$arr[3,12] = 0.06415006200365558
$arr[3,13] = 0.0641500620036556
$arr[4,0] = 14.734287
$arr[4,1] = 44.202861
$arr[4,2] = 0.4000023944294819
$arr[4,3] = 0.400002394429482
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 10.513928
$arr[4,7] = 31.541784
$arr[4,8] = 0.08295814932067838
$arr[4,9] = 0.08295814932067838
$arr[4,10] = 154.915232649336
$arr[4,11] = 1394.237093844024
$arr[4,12] = 0.03318345836570985
$arr[4,13] = 0.03318345836570986
$arr[5,0] = 14.452944
$arr[5,1] = 43.358832
$arr[5,2] = 0.3923645715978801
$arr[5,3] = 0.3923645715978802
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 68.637375
$arr[5,7] = 205.912125
$arr[5,8] = 0.5415701538216162
$arr[5,9] = 0.5415701538216162
$arr[5,10] = 992.0121371820001
$arr[5,11] = 8928.109234638
$arr[5,12] = 0.2124929413944165
$arr[5,13] = 0.2124929413944165
$arr[6,0] = 14.452944
$arr[6,1] = 43.358832
$arr[6,2] = 0.3923645715978801
$arr[6,3] = 0.3923645715978802
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 11.05007466666667
$arr[6,7] = 33.150224
$arr[6,8] = 0.08718851262838957
$arr[6,9] = 0.08718851262838957
$arr[6,10] = 159.706110353152
$arr[6,11] = 1437.354993178368
$arr[6,12] = 0.03420968340569443
$arr[6,13] = 0.03420968340569445
$arr[7,0] = 14.452944
$arr[7,1] = 43.358832
$arr[7,2] = 0.3923645715978801
$arr[7,3] = 0.3923645715978802
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 16.21089566666667
$arr[7,7] = 48.632687
$arr[7,8] = 0.1279089892319285
$arr[7,9] = 0.1279089892319285
$arr[7,10] = 234.295167260176
$arr[7,11] = 2108.656505341584
$arr[7,12] = 0.05018695576350349
$arr[7,13] = 0.0501869557635035
$arr[8,0] = 14.452944
$arr[8,1] = 43.358832
$arr[8,2] = 0.3923645715978801
$arr[8,3] = 0.3923645715978802
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 20.32546233333333
$arr[8,7] = 60.976387
$arr[8,8] = 0.1603741949973873
$arr[8,9] = 0.1603741949973873
$arr[8,10] = 293.762768877776
$arr[8,11] = 2643.864919899984
$arr[8,12] = 0.06292515231550477
$arr[8,13] = 0.06292515231550479
$arr[9,0] = 14.452944
$arr[9,1] = 43.358832
$arr[9,2] = 0.3923645715978801
$arr[9,3] = 0.3923645715978802
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 10.513928
$arr[9,7] = 31.541784
$arr[9,8] = 0.08295814932067838
$arr[9,9] = 0.08295814932067838
$arr[9,10] = 151.957212604032
$arr[9,11] = 1367.614913436288
$arr[9,12] = 0.03254983871876094
$arr[9,13] = 0.03254983871876095
$arr[10,0] = 1.259379333333333
$arr[10,1] = 3.778138
$arr[10,2] = 0.03418928576783783
$arr[10,3] = 0.03418928576783784
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 68.637375
$arr[10,7] = 205.912125
$arr[10,8] = 0.5415701538216162
$arr[10,9] = 0.5415701538216162
$arr[10,10] = 86.44049156925001
$arr[10,11] = 777.9644241232501
$arr[10,12] = 0.01851589675233913
$arr[10,13] = 0.01851589675233913
$arr[11,0] = 1.259379333333333
$arr[11,1] = 3.778138
$arr[11,2] = 0.03418928576783783
$arr[11,3] = 0.03418928576783784
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 11.05007466666667
$arr[11,7] = 33.150224
$arr[11,8] = 0.08718851262838957
$arr[11,9] = 0.08718851262838957
$arr[11,10] = 13.91623566699022
$arr[11,11] = 125.246121002912
$arr[11,12] = 0.002980912973924749
$arr[11,13] = 0.002980912973924749
$arr[12,0] = 1.259379333333333
$arr[12,1] = 3.778138
$arr[12,2] = 0.03418928576783783
$arr[12,3] = 0.03418928576783784
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 16.21089566666667
$arr[12,7] = 48.632687
$arr[12,8] = 0.1279089892319285
$arr[12,9] = 0.1279089892319285
$arr[12,10] = 20.41566697742289
$arr[12,11] = 183.741002796806
$arr[12,12] = 0.004373116985125696
$arr[12,13] = 0.004373116985125697
$arr[13,0] = 1.259379333333333
$arr[13,1] = 3.778138
$arr[13,2] = 0.03418928576783783
$arr[13,3] = 0.03418928576783784
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 20.32546233333333
$arr[13,7] = 60.976387
$arr[13,8] = 0.1603741949973873
$arr[13,9] = 0.1603741949973873
$arr[13,10] = 25.59746720304511
$arr[13,11] = 230.377204827406
$arr[13,12] = 0.005483079182552624
$arr[13,13] = 0.005483079182552626
$arr[14,0] = 1.259379333333333
$arr[14,1] = 3.778138
$arr[14,2] = 0.03418928576783783
$arr[14,3] = 0.03418928576783784
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 10.513928
$arr[14,7] = 31.541784
$arr[14,8] = 0.08295814932067838
$arr[14,9] = 0.08295814932067838
$arr[14,10] = 13.24102363535467
$arr[14,11] = 119.169212718192
$arr[14,12] = 0.002836279873895635
$arr[14,13] = 0.002836279873895636
$arr[15,0] = 4.524255666666667
$arr[15,1] = 13.572767
$arr[15,2] = 0.1228232556945456
$arr[15,3] = 0.1228232556945456
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 68.637375
$arr[15,7] = 205.912125
$arr[15,8] = 0.5415701538216162
$arr[15,9] = 0.5415701538216162
$arr[15,10] = 310.5330327888751
$arr[15,11] = 2794.797295099875
$arr[15,12] = 0.06651740947936675
$arr[15,13] = 0.06651740947936677
$arr[16,0] = 4.524255666666667
$arr[16,1] = 13.572767
$arr[16,2] = 0.1228232556945456
$arr[16,3] = 0.1228232556945456
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 11.05007466666667
$arr[16,7] = 33.150224
$arr[16,8] = 0.08718851262838957
$arr[16,9] = 0.08718851262838957
$arr[16,10] = 49.99336292775645
$arr[16,11] = 449.940266349808
$arr[16,12] = 0.01070877698018381
$arr[16,13] = 0.01070877698018381
$arr[17,0] = 4.524255666666667
$arr[17,1] = 13.572767
$arr[17,2] = 0.1228232556945456
$arr[17,3] = 0.1228232556945456
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 16.21089566666667
$arr[17,7] = 48.632687
$arr[17,8] = 0.1279089892319285
$arr[17,9] = 0.1279089892319285
$arr[17,10] = 73.3422365816588
$arr[17,11] = 660.0801292349291
$arr[17,12] = 0.01571019849006403
$arr[17,13] = 0.01571019849006404
$arr[18,0] = 4.524255666666667
$arr[18,1] = 13.572767
$arr[18,2] = 0.1228232556945456
$arr[18,3] = 0.1228232556945456
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 20.32546233333333
$arr[18,7] = 60.976387
$arr[18,8] = 0.1603741949973873
$arr[18,9] = 0.1603741949973873
$arr[18,10] = 91.95758813920322
$arr[18,11] = 827.6182932528291
$arr[18,12] = 0.01969768075897102
$arr[18,13] = 0.01969768075897102
$arr[19,0] = 4.524255666666667
$arr[19,1] = 13.572767
$arr[19,2] = 0.1228232556945456
$arr[19,3] = 0.1228232556945456
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 10.513928
$arr[19,7] = 31.541784
$arr[19,8] = 0.08295814932067838
$arr[19,9] = 0.08295814932067838
$arr[19,10] = 47.56769833292534
$arr[19,11] = 428.109284996328
$arr[19,12] = 0.01018918998595997
$arr[19,13] = 0.01018918998595998
$arr[20,0] = 1.864631
$arr[20,1] = 5.593893
$arr[20,2] = 0.05062049251025444
$arr[20,3] = 0.05062049251025445
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 68.637375
$arr[20,7] = 205.912125
$arr[20,8] = 0.5415701538216162
$arr[20,9] = 0.5415701538216162
$arr[20,10] = 127.983377183625
$arr[20,11] = 1151.850394652625
$arr[20,12] = 0.02741454791530447
$arr[20,13] = 0.02741454791530447
$arr[21,0] = 1.864631
$arr[21,1] = 5.593893
$arr[21,2] = 0.05062049251025444
$arr[21,3] = 0.05062049251025445
$arr[21,4] = 3
$arr[21,5] = 1
$arr[21,6] = 11.05007466666667
$arr[21,7] = 33.150224
$arr[21,8] = 0.08718851262838957
$arr[21,9] = 0.08718851262838957
$arr[21,10] = 20.60431177578133
$arr[21,11] = 185.438805982032
$arr[21,12] = 0.004413525450485619
$arr[21,13] = 0.00441352545048562
$arr[22,0] = 1.864631
$arr[22,1] = 5.593893
$arr[22,2] = 0.05062049251025444
$arr[22,3] = 0.05062049251025445
$arr[22,4] = 3
$arr[22,5] = 1
$arr[22,6] = 16.21089566666667
$arr[22,7] = 48.632687
$arr[22,8] = 0.1279089892319285
$arr[22,9] = 0.1279089892319285
$arr[22,10] = 30.22733859783234
$arr[22,11] = 272.046047380491
$arr[22,12] = 0.006474816031409053
$arr[22,13] = 0.006474816031409054
$arr[23,0] = 1.864631
$arr[23,1] = 5.593893
$arr[23,2] = 0.05062049251025444
$arr[23,3] = 0.05062049251025445
$arr[23,4] = 3
$arr[23,5] = 1
$arr[23,6] = 20.32546233333333
$arr[23,7] = 60.976387
$arr[23,8] = 0.1603741949973873
$arr[23,9] = 0.1603741949973873
$arr[23,10] = 37.89948715606567
$arr[23,11] = 341.095384404591
$arr[23,12] = 0.008118220736703331
$arr[23,13] = 0.008118220736703331
$arr[24,0] = 1.864631
$arr[24,1] = 5.593893
$arr[24,2] = 0.05062049251025444
$arr[24,3] = 0.05062049251025445
$arr[24,4] = 3
$arr[24,5] = 1
$arr[24,6] = 10.513928
$arr[24,7] = 31.541784
$arr[24,8] = 0.08295814932067838
$arr[24,9] = 0.08295814932067838
$arr[24,10] = 19.604596080568
$arr[24,11] = 176.441364725112
$arr[24,12] = 0.004199382376351969
$arr[24,13] = 0.00419938237635197

$ws.Range("G2:T26").Value2 = $arr
